$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("TestCases")

# Row 3 (H3): AppMinimized=Homescreen -> com.symbol.enterprisebrowser
$ws.Range("H3").Value = "validate1`n{`nvalidate_PageTitle=Manual specs`n};`nvalidate2`n{`nvalidate_PageTitle=Application JS Test`n};`nvalidate3`n{`nvalidate_Text_Exists=VT281-0007`n};`nvalidate5`n{`nvalidate_PageTitle=Manual specs`n};`nvalidate6`n{`nvalidate_PageTitle=Application JS Test`n};`nvalidate7`n{`nvalidate_Page=settingspage_xpath,Client ID`n};`nvalidate8`n{`nvalidate_AppMinimized=com.symbol.enterprisebrowser`n};"

# Row 6 (G6): wait(20) -> wait(30)
$ws.Range("G6").Value = "wait(3);`nvalidate1;`nlink_Click(Application_test_link);`nvalidate2;`nSelectTestToRun(VT281_0020_string);`nClickRunTest(runtest_top_xpath);`nvalidate3;`nClickRunTest(runtest_bottom_xpath);`nwait(30);`nvalidate4;"

# Row 10 (H10): AppMinimized=Homescreen -> com.symbol.enterprisebrowser
$ws.Range("H10").Value = "validate1`n{`nvalidate_PageTitle=Manual specs`n};`nvalidate2`n{`nvalidate_PageTitle=Application JS Test`n};`nvalidate3`n{`nvalidate_Text_Exists=VT281-0031`n};`nvalidate4`n{`nvalidate_AppMinimized=com.symbol.enterprisebrowser`n};`n"

# Row 11 (H11): AppMinimized=Appsscreen -> com.symbol.enterprisebrowser
$ws.Range("H11").Value = "validate1`n{`nvalidate_PageTitle=Manual specs`n};`nvalidate2`n{`nvalidate_PageTitle=Application JS Test`n};`nvalidate3`n{`nvalidate_Text_Exists=VT281-0032`n};`nvalidate4`n{`nvalidate_AppMinimized=com.symbol.enterprisebrowser`n};`n"

# Row 12 (H12): AppMinimized=Homescreen -> com.symbol.enterprisebrowser
$ws.Range("H12").Value = "validate1`n{`nvalidate_PageTitle=Manual specs`n};`nvalidate2`n{`nvalidate_PageTitle=Application JS Test`n};`nvalidate3`n{`nvalidate_Text_Exists=VT281-0069`n};`nvalidate4`n{`nvalidate_Page=settingspage_xpath,Client ID`n};`nvalidate5`n{`nvalidate_PageTitle=Application JS Test`n};`nvalidate6`n{`nvalidate_Logdisplayed=LOG_TEST`nvalidate_Logdisplayed=menuCallback`n};`nvalidate7`n{`nvalidate_AppMinimized=com.symbol.enterprisebrowser`n};`n`n"

# Row 13 (H13): AppMinimized=Homescreen -> com.symbol.enterprisebrowser
$ws.Range("H13").Value = "validate1`n{`nvalidate_PageTitle=Manual specs`n};`nvalidate2`n{`nvalidate_PageTitle=Application JS Test`n};`nvalidate3`n{`nvalidate_Text_Exists=VT281-0070`n};`nvalidate4`n{`nvalidate_PageTitle=Application JS Test`n};`nvalidate5`n{`nvalidate_AppMinimized=com.symbol.enterprisebrowser`n};`n`n"

# Row 16 (H16): AppMinimized=AppsScreen -> com.symbol.enterprisebrowser
$ws.Range("H16").Value = "validate1`n{`nvalidate_PageTitle=Manual specs`n};`nvalidate2`n{`nvalidate_PageTitle=Application JS Test`n};`nvalidate3`n{`nvalidate_Text_Exists=VT281-0077`n};`nvalidate4`n{`nvalidate_AppMinimized=com.symbol.enterprisebrowser`n};`nvalidate5`n{`nvalidate_Result=UIDestroyed`nvalidate_Result=Deactivated`nvalidate_Result=Activated`nvalidate_doesNotContain=ScreenOff`nvalidate_doesNotContain=ScreenOn`n};"

# Row 18 (H18): AppMinimized=AppsScreen -> com.symbol.enterprisebrowser
$ws.Range("H18").Value = "validate1`n{`nvalidate_PageTitle=Manual specs`n};`nvalidate2`n{`nvalidate_PageTitle=Application JS Test`n};`nvalidate3`n{`nvalidate_Text_Exists=VT281-0079`n};`nvalidate4`n{`nvalidate_AppMinimized=com.symbol.enterprisebrowser`n};`nvalidate5`n{`nvalidate_PageTitle=Application JS Test`n};`nvalidate6`n{`nvalidate_Result=UIDestroyed`nvalidate_Result=Deactivated`nvalidate_Result=ScreenOff`nvalidate_Result=Activated`nvalidate_Result=ScreenOn`n};"

$ws.Range("A2").Select()
